$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row that only held "2346890 - Eliane Corrêa Pedrozo" in B13/C13 (with no
# A-column label) is removed entirely; everything below shifts up by one row.
$ws.Rows(13).Delete()

# --- Content updates on top of the shifted layout ---

# Objetivos: row now holds the docente info instead of the long paragraph.
$ws.Range("B10").Value = "2346890 - Eliane Corrêa Pedrozo"
$ws.Range("C10").Value = "2346890 - Eliane Corrêa Pedrozo"

# Programa resumido: (former row 14) now just says "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Programa: (former row 16) gains a date value that wasn't there before -
# copy formatting from an existing B/C pair first so the new cells pick up
# the correct (wrap-text / red-text) column styles.
$ws.Range("B10:C10").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("B15").Value = "15/07/2015"
$ws.Range("C15").Value = "15/07/2015"

# Método: (former row 19) now shows the docente info.
$ws.Range("B18").Value = "2346890 - Eliane Corrêa Pedrozo"
$ws.Range("C18").Value = "2346890 - Eliane Corrêa Pedrozo"

# Critério: (former row 20) now holds the seminar/evaluation method text.
$seminario = "Seminário em grupo sobre um estudo de caso apresentado pelos alunos.`nResolução de um exercício individual após cada aula sobre o tema abordado, com consulta. `nProva escrita."
$ws.Range("B19").Value = $seminario
$ws.Range("C19").Value = $seminario

# Norma de recuperação: (former row 21) now holds the "Média Final" formula text.
$mediaFinal = "Média Final = 0,4 x Nota da Prova + 0,2 x Nota dos exercícios + 0,4 x Nota do Seminário`n`nMédia Final Mínima para Aprovação = 5,0"
$ws.Range("B20").Value = $mediaFinal
$ws.Range("C20").Value = $mediaFinal

# Bibliografia: (former row 22) now holds the "Nota Final" formula text.
$notaFinal = "Nota Final = (Prova Escrita + Média final)/2`n`nNota Final Mínima para Aprovação = 5,0"
$ws.Range("B21").Value = $notaFinal
$ws.Range("C21").Value = $notaFinal
